$wb = $excel.ActiveWorkbook

# --- 1. Copy the MeetingResume sheet (keeps formatting/col widths/styles identical) ---
$meetingResume = $wb.Sheets.Item("MeetingResume")
$lastSheet = $wb.Sheets.Item($wb.Sheets.Count)
$meetingResume.Copy([System.Reflection.Missing]::Value, $lastSheet)
$newSheet = $wb.Sheets.Item($wb.Sheets.Count)
$newSheet.Name = "InvokeviaTxt"

# --- 2. Update the body/url content on the copied sheet ---
$url = "/cvi/dm/api/v1/invoke/text/json?intent=true&skill=true"
$bodyToronto = "{`n  ""text"": ""starte das toronto meeting""`n}"
$bodyMesse = "{`n  ""text"": ""starte das Messe meeting""`n}"

$newSheet.Range("A2").Value = $url
$newSheet.Range("E2").Value = $bodyToronto

$newSheet.Range("A3").Value = $url
$newSheet.Range("E3").Value = $bodyMesse

# --- 3. Row heights for the shorter JSON text ---
$newSheet.Rows.Item(2).RowHeight = 43.5
$newSheet.Rows.Item(3).RowHeight = 43.5

# --- 4. Restore whole-sheet selection on MeetingResume (it lost the active-tab status) ---
$meetingResume.Activate()
$meetingResume.Cells.Select()

# --- 5. Select B3 on the new sheet and make it the active tab ---
$newSheet.Activate()
$newSheet.Range("B3").Select()
